$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 2024 update added two new study rows into what was an unused blank-row
# gap (rows 66-69) below row 65. One of those now-unneeded blank rows is
# removed so the remaining rows below (the stray whitespace row + the trailing
# blank rows) shift up by one, matching the final layout.
$ws.Rows.Item(66).Delete()

# Populate the new content in the order it was originally typed/pasted so the
# shared-string table ends up in the same sequence.
$ws.Cells.Item(67,1).Value2 = "Multi-objective optimization and integrated indicator-
driven two-stage project recommendation in time-
dependent software ecosystem"
$ws.Cells.Item(66,1).Value2 = "CodeCompass: NLP-Driven Navigation to Optimal Repositories"
$ws.Cells.Item(66,16).Value2 = "Sheetal Phatangare, Aakash Matkar,Akshay Jadhav
"
$ws.Cells.Item(66,18).Value2 = "2024 4th International Conference on Pervasive Computing and Social Networking (ICPCSN)"
$ws.Cells.Item(67,18).Value2 = "Information and Software Technology"
$ws.Cells.Item(67,16).Value2 = "Xin Shen a, Xiangjuan Yao b, Dunwei Gong c,∗, Huijie Tu"
$ws.Cells.Item(66,12).Value2 = "repository metadata, programming languages used,
commit history, issue tracking, and community engagement
metrics."
$ws.Cells.Item(66,4).Value2 = "The system capitalizes on various NLP
methodologies, including lemmat ization, stemming,
tokenizat ion, and the innovative Bidirectional Encoder
Representations from Transformers (BERT)."
$ws.Cells.Item(66,5).Value2 = "vector based"
$ws.Cells.Item(66,6).Value2 = "By leveraging sophisticated NLP
techniques and semantic analysis, it provides developers
with tailored recommendations, minimizing the effort
required to discover suitable projects to contribute to within
the vast GitHub ecosystem. Additionally, the interactive
web applicat ion fosters a user-friendly environment,
encouraging collaboration and facilitating informed
decision-making."
$ws.Cells.Item(66,7).Value2 = "Additionally, integrating large language models like GPT
(Generative Pre-trained Transformer) could further enhance
accuracy. Large language models have shown remarkable
performance in various NLP tasks due to their ability to
capture complex linguistic patterns and semantics. By finetuning
a pre-trained GPT model on specific GitHub issue
data, we can create a more sophisticated recommendation
system that understands the nuances of the text and provides
more accurate suggestions."
$ws.Cells.Item(67,5).Value2 = "hybrid"
$ws.Cells.Item(67,2).Value2 = "behaviours, social connections"
$ws.Cells.Item(67,12).Value2 = "programming language, timestamps of their contributions,"
$ws.Cells.Item(67,7).Value2 = "However, this paper
mines the preferences of these developers and their communities from
the programming language and profession perspective.
some
developers and projects often lack historical information, which brings
greater challenges to project recommendation. Therefore, how to tackle
the problem will be another topic to be further studied."
$ws.Cells.Item(67,6).Value2 = "the proposed method obtains
better success rate and efficiency of recommendation compared with
comparison ones."
$ws.Cells.Item(67,11).Value2 = "comments, PR"
$ws.Cells.Item(67,4).Value2 = "improved NSGA-II algorithm,
The method proposed in the paper for project recommendation in a time-dependent software ecosystem involves a two-stage framework. 
1. Initial Setup:
- Start with the current state of developers and their projects at the initial time interval (t = 0).
2. Change Detection:
- For each subsequent time interval (t = 1 to T), monitor and detect changes in developer communities and projects.
3. Recommendation Logic:
- If neither developers nor projects have changed, retain the previous recommendations.
- If only the projects have changed, use the existing developer community to recommend new projects.
- If both the developer community and projects have changed, re-divide the community and update the project recommendations accordingly.
4. Output Recommendations:
- Upon completing the iterations over T time intervals, output the recommended developers and their associated projects."
$ws.Cells.Item(66,2).Value2 = "programming preference"
$ws.Cells.Item(66,9).Value2 = "Skill preference "
$ws.Cells.Item(67,9).Value2 = "Social and behavioral preferences "
$ws.Cells.Item(66,11).Value2 = "comments "

# Remaining cells reuse values already present elsewhere in the sheet.
$ws.Cells.Item(66,3).Value2 = "project features"
$ws.Cells.Item(66,8).Value2 = "repository to work on"
$ws.Cells.Item(66,10).Value2 = "no suggestions"
$ws.Cells.Item(66,13).Value2 = "Github"
$ws.Cells.Item(66,14).Value2 = "not mentioned"
$ws.Cells.Item(67,3).Value2 = "project features"
$ws.Cells.Item(67,8).Value2 = "repository to work on"
$ws.Cells.Item(67,10).Value2 = "no suggestions"
$ws.Cells.Item(67,13).Value2 = "Github"
$ws.Cells.Item(67,14).Value2 = "not mentioned"

# Year + Q column (plain numbers, no shared string).
$ws.Cells.Item(66,17).Value2 = 2024
$ws.Cells.Item(67,17).Value2 = 2024

# Wrap text matches the rest of the data rows.
$row66WrapCols = 1,2,4,5,6,7,8,9,10,11,12,13,14,16
$row67WrapCols = 1,2,4,5,6,7,8,9,10,11,12,13,14
foreach ($col in $row66WrapCols) { $ws.Cells.Item(66, $col).WrapText = $true }
foreach ($col in $row67WrapCols) { $ws.Cells.Item(67, $col).WrapText = $true }

# Auto row heights for the newly entered, word-wrapped content.
$ws.Rows.Item(66).RowHeight = 246.5
$ws.Rows.Item(67).RowHeight = 409.5

# Leave the selection where the editor finished up.
$ws.Range("K67").Select()
